# "update: csv and xlsx file for the student id to int"
#
# The sheet's column A held student IDs stored as 5-digit numbers
# (70001..70100). The commit converts them to plain sequential integers
# (1..100), leaving columns B-D (the TRUE_NEW/SKIP/... labels) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 1; $row -le 100; $row++) {
    $ws.Cells.Item($row, 1).Value = $row
}

# The diff also shows the saved selection moved from E3 to D6 - mirror that.
$ws.Range("D6").Select()
